$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (bold font, borders, centered/top
# aligned) onto the two new header cells so they reuse the same style as the
# rest of the header row instead of creating brand new style entries.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
